$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price column cells we are about to rewrite,
# so values like "598.17" are not auto-converted to numbers by Excel.
$priceCells = @("D2","D3","D5","D6","D9","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D40","D41","D43","D45","D46","D47","D48","D49","D51")
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "72.127.00"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.687.58"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "598.17"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "174.23"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "2.687.52"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "5.00"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "3.177.73"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000184"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "71.904.79"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "26.22"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "2.686.44"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "12.26"
$ws.Range("E19").Value = "  +6.25%  "
$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "371.11"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "72.38"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D26").Value = "4.34"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").Value = "9.76"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "2.821.45"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "0.0₃0963"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "8.07"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "500.18"
$ws.Range("E32").Value = "  -9.08%  "
$ws.Range("D33").Value = "1.30"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "164.16"
$ws.Range("D37").Value = "19.61"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "19.10"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").Value = "0.108"
$ws.Range("E40").Value = "  -6.17%  "
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "5.02"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").Value = "157.36"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "39.41"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "0.567"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").Value = "0.0762"
$ws.Range("E51").Value = "  +0.82%  "

# Restore the default (General/Normal) cell style now that the text is stored,
# so formatting matches the original workbook.
foreach ($addr in $priceCells) { $ws.Range($addr).Style = "Normal" }
